$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add 4 new rows (27-30) by copying the formatting of row 26 ---
$ws.Range("A26:K26").Copy($ws.Range("A27:K27"))
$ws.Range("A26:K26").Copy($ws.Range("A28:K28"))
$ws.Range("A26:K26").Copy($ws.Range("A29:K29"))
$ws.Range("A26:K26").Copy($ws.Range("A30:K30"))

# --- 2. Update the shared B..I values (same for every data row 2-30) ---
$origin = "BA11 5LB"
$destination = "BA11 5AP"
$startAddress = "81 Knights Maltings, Frome, Frome, BA11 5LB, United Kingdom"
$endAddress = "55 Tower View, Frome, Frome, BA11 5AP, United Kingdom"
$distanceText = 3.0501
$distanceValue = 3050.1
$durationText = 8.711666666666668
$durationValue = 522.7

for ($row = 2; $row -le 30; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 2
    $ws.Cells.Item($row, 2).Value = $origin
    $ws.Cells.Item($row, 3).Value = $destination
    $ws.Cells.Item($row, 4).Value = $startAddress
    $ws.Cells.Item($row, 5).Value = $endAddress
    $ws.Cells.Item($row, 6).Value = $distanceText
    $ws.Cells.Item($row, 7).Value = $distanceValue
    $ws.Cells.Item($row, 8).Value = $durationText
    $ws.Cells.Item($row, 9).Value = $durationValue
}

# --- 3. Update the per-row Lat/Lng (columns J/K) ---
$ws.Cells.Item(2, 10).Value = 51.22234
$ws.Cells.Item(2, 11).Value = -2.31109
$ws.Cells.Item(3, 10).Value = 51.22237
$ws.Cells.Item(3, 11).Value = -2.3107
$ws.Cells.Item(4, 10).Value = 51.22273
$ws.Cells.Item(4, 11).Value = -2.31064
$ws.Cells.Item(5, 10).Value = 51.22283
$ws.Cells.Item(5, 11).Value = -2.31005
$ws.Cells.Item(6, 10).Value = 51.22298
$ws.Cells.Item(6, 11).Value = -2.30982
$ws.Cells.Item(7, 10).Value = 51.22374
$ws.Cells.Item(7, 11).Value = -2.30909
$ws.Cells.Item(8, 10).Value = 51.22498
$ws.Cells.Item(8, 11).Value = -2.30754
$ws.Cells.Item(9, 10).Value = 51.22534
$ws.Cells.Item(9, 11).Value = -2.30686
$ws.Cells.Item(10, 10).Value = 51.22581
$ws.Cells.Item(10, 11).Value = -2.3054
$ws.Cells.Item(11, 10).Value = 51.22681
$ws.Cells.Item(11, 11).Value = -2.30373
$ws.Cells.Item(12, 10).Value = 51.22708
$ws.Cells.Item(12, 11).Value = -2.30363
$ws.Cells.Item(13, 10).Value = 51.22726
$ws.Cells.Item(13, 11).Value = -2.30377
$ws.Cells.Item(14, 10).Value = 51.22884
$ws.Cells.Item(14, 11).Value = -2.3063
$ws.Cells.Item(15, 10).Value = 51.22893
$ws.Cells.Item(15, 11).Value = -2.30699
$ws.Cells.Item(16, 10).Value = 51.22876
$ws.Cells.Item(16, 11).Value = -2.30829
$ws.Cells.Item(17, 10).Value = 51.22791
$ws.Cells.Item(17, 11).Value = -2.31099
$ws.Cells.Item(18, 10).Value = 51.22768
$ws.Cells.Item(18, 11).Value = -2.31325
$ws.Cells.Item(19, 10).Value = 51.22642
$ws.Cells.Item(19, 11).Value = -2.31437
$ws.Cells.Item(20, 10).Value = 51.22582
$ws.Cells.Item(20, 11).Value = -2.31544
$ws.Cells.Item(21, 10).Value = 51.22519
$ws.Cells.Item(21, 11).Value = -2.31769
$ws.Cells.Item(22, 10).Value = 51.22484
$ws.Cells.Item(22, 11).Value = -2.32133
$ws.Cells.Item(23, 10).Value = 51.22421
$ws.Cells.Item(23, 11).Value = -2.32144
$ws.Cells.Item(24, 10).Value = 51.22301
$ws.Cells.Item(24, 11).Value = -2.32124
$ws.Cells.Item(25, 10).Value = 51.22183
$ws.Cells.Item(25, 11).Value = -2.32125
$ws.Cells.Item(26, 10).Value = 51.22048
$ws.Cells.Item(26, 11).Value = -2.32074
$ws.Cells.Item(27, 10).Value = 51.22034
$ws.Cells.Item(27, 11).Value = -2.31956
$ws.Cells.Item(28, 10).Value = 51.21988
$ws.Cells.Item(28, 11).Value = -2.31827
$ws.Cells.Item(29, 10).Value = 51.22045
$ws.Cells.Item(29, 11).Value = -2.31728
$ws.Cells.Item(30, 10).Value = 51.22032
$ws.Cells.Item(30, 11).Value = -2.31717

Write-Host "Done. Used range:" $ws.UsedRange.Address()
